$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.334.96'
$ws.Cells.Item(2, 5).Value = '  +3.73%  '
$ws.Cells.Item(3, 4).Value = '1.721.08'
$ws.Cells.Item(3, 5).Value = '  +3.50%  '
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.9990'
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '239.83'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.70%  '
$ws.Cells.Item(6, 5).Value = '  -0.04%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.4721'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.28%  '
$ws.Cells.Item(8, 5).Value = '  +0.78%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.06222'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.17%  '
$ws.Cells.Item(10, 4).Value = '1.715.30'
$ws.Cells.Item(10, 5).Value = '  +3.14%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.07074'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.05%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '15.26'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +3.51%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.5924'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.22%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '4.411'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.65%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '76.36'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.64%  '
$ws.Cells.Item(16, 5).Value = '  -0.01%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.03%  '
$ws.Cells.Item(18, 4).Value = '26.328.49'
$ws.Cells.Item(18, 5).Value = '  +3.71%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.000006802'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.58%  '
$ws.Cells.Item(20, 5).Value = '  +1.71%  '
$ws.Cells.Item(21, 4).Value = '1.936.52'
$ws.Cells.Item(21, 5).Value = '  +3.40%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '4.561'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.65%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '8.788'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.71%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '5.339'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '134.92'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.11%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '15.19'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.94%  '
$ws.Cells.Item(27, 5).Value = '  +0.20%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '1.763'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.52%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '106.74'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +2.42%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '4.024'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.84%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '3.694'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.08%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.07732'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.12%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '0.04448'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.71%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '2.612'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.23%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.9745'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.27%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.6205'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.48%  '
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.9275'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +8.66%  '
$ws.Cells.Item(38, 2).Value = 'Quant'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '114.53'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +16.69%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.415'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -7.61%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.05%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '1.906'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +4.19%  '
$ws.Cells.Item(42, 5).Value = '  -2.21%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '5.293'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +13.55%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.3821'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.50%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.1159'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +4.51%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '6.256'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.83%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.05295'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.83%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '30.64'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.79%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '7.685'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +4.89%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.3391'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.47%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.220'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +1.52%  '
